$p = $ppt.ActivePresentation

# The commit removes the "Picture Placeholder" slide (SlideID 272, which
# sat at position 16, right before the closing "THANK YOU" slide,
# SlideID 270) from the deck. Locate it by its stable SlideID and delete
# it; PowerPoint then renumbers the sldIdLst/relationship bookkeeping and
# the trailing "THANK YOU" slide (270) slides up to become the new last
# slide.
for ($i = $p.Slides.Count; $i -ge 1; $i--) {
    $s = $p.Slides.Item($i)
    if ($s.SlideID -eq 272) {
        $s.Delete()
    }
}
